$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.339138636047295
$ws.Cells.Item(2, 3).Value2 = 0.6822874540790167
$ws.Cells.Item(2, 4).Value2 = 0.6876767610321508
$ws.Cells.Item(2, 5).Value2 = 0.2797642199864896
$ws.Cells.Item(2, 7).Value2 = 0.002631925694092923
$ws.Cells.Item(2, 9).Value2 = 2.315527991889027
$ws.Cells.Item(2, 10).Value2 = 0.1450022493354552
$ws.Cells.Item(2, 11).Value2 = 2.019307990800598
$ws.Cells.Item(2, 14).Value2 = 3.92459706957834
$ws.Cells.Item(3, 2).Value2 = 1.295990837717596
$ws.Cells.Item(3, 3).Value2 = 0.6656458095005746
$ws.Cells.Item(3, 4).Value2 = 0.67903188140059
$ws.Cells.Item(3, 5).Value2 = 0.2754301794278007
$ws.Cells.Item(3, 7).Value2 = 0.002637331869385179
$ws.Cells.Item(3, 9).Value2 = 2.297489746990578
$ws.Cells.Item(3, 10).Value2 = 0.1420633201189716
$ws.Cells.Item(3, 11).Value2 = 1.959706342893128
$ws.Cells.Item(3, 14).Value2 = 3.912949913859421
$ws.Cells.Item(4, 2).Value2 = 1.270379787769343
$ws.Cells.Item(4, 3).Value2 = 0.6558268404656644
$ws.Cells.Item(4, 4).Value2 = 0.6740936387296301
$ws.Cells.Item(4, 5).Value2 = 0.2729267394352632
$ws.Cells.Item(4, 7).Value2 = 0.002640824128888543
$ws.Cells.Item(4, 9).Value2 = 2.287210084993347
$ws.Cells.Item(4, 10).Value2 = 0.1403457294324113
$ws.Cells.Item(4, 11).Value2 = 1.924393212673806
$ws.Cells.Item(4, 14).Value2 = 3.906529951568572
$ws.Cells.Item(5, 2).Value2 = 1.260164107528794
$ws.Cells.Item(5, 3).Value2 = 0.6519254111217663
$ws.Cells.Item(5, 4).Value2 = 0.672174030531238
$ws.Cells.Item(5, 5).Value2 = 0.2719460968690868
$ws.Cells.Item(5, 7).Value2 = 0.00264229086948925
$ws.Cells.Item(5, 9).Value2 = 2.283220633502921
$ws.Cells.Item(5, 10).Value2 = 0.1396675628127468
$ws.Cells.Item(5, 11).Value2 = 1.910324124767698
$ws.Cells.Item(5, 14).Value2 = 3.904097054783151
$ws.Cells.Item(6, 2).Value2 = 1.258481130554941
$ws.Cells.Item(6, 3).Value2 = 0.651283600748684
$ws.Cells.Item(6, 4).Value2 = 0.6718608778506905
$ws.Cells.Item(6, 5).Value2 = 0.2717856459559655
$ws.Cells.Item(6, 7).Value2 = 0.002642537059569579
$ws.Cells.Item(6, 9).Value2 = 2.282570227559006
$ws.Cells.Item(6, 10).Value2 = 0.1395562656874674
$ws.Cells.Item(6, 11).Value2 = 1.908007334577974
$ws.Cells.Item(6, 14).Value2 = 3.903704125923298
$ws.Cells.Item(7, 2).Value2 = 1.270241121542227
$ws.Cells.Item(7, 3).Value2 = 0.6557738206303156
$ws.Cells.Item(7, 4).Value2 = 0.674067374889546
$ws.Cells.Item(7, 5).Value2 = 0.2729133542484519
$ws.Cells.Item(7, 7).Value2 = 0.002640843732953352
$ws.Cells.Item(7, 9).Value2 = 2.287155474439103
$ws.Cells.Item(7, 10).Value2 = 0.1403364954317823
$ws.Cells.Item(7, 11).Value2 = 1.924202172687131
$ws.Cells.Item(7, 14).Value2 = 3.906496399420007
$ws.Cells.Item(8, 2).Value2 = 1.324077784991744
$ws.Cells.Item(8, 3).Value2 = 0.6764662973993438
$ws.Cells.Item(8, 4).Value2 = 0.6846191326024496
$ws.Cells.Item(8, 5).Value2 = 0.27823704167767
$ws.Cells.Item(8, 7).Value2 = 0.002633753961345079
$ws.Cells.Item(8, 9).Value2 = 2.309142880572622
$ws.Cells.Item(8, 10).Value2 = 0.1439708083603151
$ws.Cells.Item(8, 11).Value2 = 1.998490401031574
$ws.Cells.Item(8, 14).Value2 = 3.920428928475019
$ws.Cells.Item(9, 2).Value2 = 1.436689252943609
$ws.Cells.Item(9, 3).Value2 = 0.7202349373309005
$ws.Cells.Item(9, 4).Value2 = 0.7082565362940443
$ws.Cells.Item(9, 5).Value2 = 0.2899342849927891
$ws.Cells.Item(9, 7).Value2 = 0.00262121539946891
$ws.Cells.Item(9, 9).Value2 = 2.35860392767907
$ws.Cells.Item(9, 10).Value2 = 0.1517921466171117
$ws.Cells.Item(9, 11).Value2 = 2.154411013285653
$ws.Cells.Item(9, 14).Value2 = 3.953586157510813
$ws.Cells.Item(10, 2).Value2 = 1.523783349971779
$ws.Cells.Item(10, 3).Value2 = 0.7543752349287161
$ws.Cells.Item(10, 4).Value2 = 0.7274378169623503
$ws.Cells.Item(10, 5).Value2 = 0.2993051570506964
$ws.Cells.Item(10, 7).Value2 = 0.002612825314130422
$ws.Cells.Item(10, 9).Value2 = 2.398856522454921
$ws.Cells.Item(10, 10).Value2 = 0.1579693135705611
$ws.Cells.Item(10, 11).Value2 = 2.27531617554655
$ws.Cells.Item(10, 14).Value2 = 3.981553241874025
$ws.Cells.Item(11, 2).Value2 = 1.564366780355215
$ws.Cells.Item(11, 3).Value2 = 0.770345925819413
$ws.Cells.Item(11, 4).Value2 = 0.736562415103748
$ws.Cells.Item(11, 5).Value2 = 0.3037392640753183
$ws.Cells.Item(11, 7).Value2 = 0.00260918483763939
$ws.Cells.Item(11, 9).Value2 = 2.418028815376886
$ws.Cells.Item(11, 10).Value2 = 0.1608747102681605
$ws.Cells.Item(11, 11).Value2 = 2.331722432041772
$ws.Cells.Item(11, 14).Value2 = 3.995070297751084
$ws.Cells.Item(12, 2).Value2 = 1.579874355068853
$ws.Cells.Item(12, 3).Value2 = 0.7764575304545929
$ws.Cells.Item(12, 4).Value2 = 0.740075345139303
$ws.Cells.Item(12, 5).Value2 = 0.305443140508558
$ws.Cells.Item(12, 7).Value2 = 0.002607831462937679
$ws.Cells.Item(12, 9).Value2 = 2.425413459818159
$ws.Cells.Item(12, 10).Value2 = 0.1619887469788353
$ws.Cells.Item(12, 11).Value2 = 2.353285901171262
$ws.Cells.Item(12, 14).Value2 = 4.000303953595846
$ws.Cells.Item(13, 2).Value2 = 1.576528301968551
$ws.Cells.Item(13, 3).Value2 = 0.7751384381504636
$ws.Cells.Item(13, 4).Value2 = 0.7393162046481052
$ws.Cells.Item(13, 5).Value2 = 0.3050750758596195
$ws.Cells.Item(13, 7).Value2 = 0.002608121818330034
$ws.Cells.Item(13, 9).Value2 = 2.423817495827805
$ws.Cells.Item(13, 10).Value2 = 0.1617482025230288
$ws.Cells.Item(13, 11).Value2 = 2.348632741231938
$ws.Cells.Item(13, 14).Value2 = 3.999171663406912
$ws.Cells.Item(14, 2).Value2 = 1.565639796598816
$ws.Cells.Item(14, 3).Value2 = 0.7708474482705014
$ws.Cells.Item(14, 4).Value2 = 0.7368502691174115
$ws.Cells.Item(14, 5).Value2 = 0.3038789455837474
$ws.Cells.Item(14, 7).Value2 = 0.002609072990621408
$ws.Cells.Item(14, 9).Value2 = 2.418633855168153
$ws.Cells.Item(14, 10).Value2 = 0.1609660849968009
$ws.Cells.Item(14, 11).Value2 = 2.333492381702001
$ws.Cells.Item(14, 14).Value2 = 3.99549856377061
$ws.Cells.Item(15, 2).Value2 = 1.558988468279438
$ws.Cells.Item(15, 3).Value2 = 0.7682274252733237
$ws.Cells.Item(15, 4).Value2 = 0.735347327513665
$ws.Cells.Item(15, 5).Value2 = 0.3031495123985266
$ws.Cells.Item(15, 7).Value2 = 0.002609658887837263
$ws.Cells.Item(15, 9).Value2 = 2.415474961185268
$ws.Cells.Item(15, 10).Value2 = 0.1604888193383687
$ws.Cells.Item(15, 11).Value2 = 2.324245039259324
$ws.Cells.Item(15, 14).Value2 = 3.993263688703564
$ws.Cells.Item(16, 2).Value2 = 1.521150600027624
$ws.Cells.Item(16, 3).Value2 = 0.7533404166544528
$ws.Cells.Item(16, 4).Value2 = 0.7268495540768924
$ws.Cells.Item(16, 5).Value2 = 0.2990188365273951
$ws.Cells.Item(16, 7).Value2 = 0.002613066761211547
$ws.Cells.Item(16, 9).Value2 = 2.397620954713631
$ws.Cells.Item(16, 10).Value2 = 0.1577813676889406
$ws.Cells.Item(16, 11).Value2 = 2.271658309744964
$ws.Cells.Item(16, 14).Value2 = 3.980685926763869
$ws.Cells.Item(17, 2).Value2 = 1.498185804464072
$ws.Cells.Item(17, 3).Value2 = 0.7443208371111609
$ws.Cells.Item(17, 4).Value2 = 0.7217388167284753
$ws.Cells.Item(17, 5).Value2 = 0.2965287681145838
$ws.Cells.Item(17, 7).Value2 = 0.002615202410843579
$ws.Cells.Item(17, 9).Value2 = 2.386889146551908
$ws.Cells.Item(17, 10).Value2 = 0.1561449396579064
$ws.Cells.Item(17, 11).Value2 = 2.239759177903807
$ws.Cells.Item(17, 14).Value2 = 3.973173956525955
$ws.Cells.Item(18, 2).Value2 = 1.485067731722268
$ws.Cells.Item(18, 3).Value2 = 0.7391744031594953
$ws.Cells.Item(18, 4).Value2 = 0.7188367886847686
$ws.Cells.Item(18, 5).Value2 = 0.2951126562841324
$ws.Cells.Item(18, 7).Value2 = 0.002616447374943398
$ws.Cells.Item(18, 9).Value2 = 2.38079751194671
$ws.Cells.Item(18, 10).Value2 = 0.1552126796351274
$ws.Cells.Item(18, 11).Value2 = 2.221543879628427
$ws.Cells.Item(18, 14).Value2 = 3.968928023388315
$ws.Cells.Item(19, 2).Value2 = 1.480641723219776
$ws.Cells.Item(19, 3).Value2 = 0.7374390007360603
$ws.Cells.Item(19, 4).Value2 = 0.7178606498778493
$ws.Cells.Item(19, 5).Value2 = 0.2946359474652738
$ws.Cells.Item(19, 7).Value2 = 0.002616871752971317
$ws.Cells.Item(19, 9).Value2 = 2.378748883620801
$ws.Cells.Item(19, 10).Value2 = 0.1548985696697116
$ws.Cells.Item(19, 11).Value2 = 2.215399160247102
$ws.Cells.Item(19, 14).Value2 = 3.967503239228279
$ws.Cells.Item(20, 2).Value2 = 1.500621055653312
$ws.Cells.Item(20, 3).Value2 = 0.745276699024771
$ws.Cells.Item(20, 4).Value2 = 0.7222789769571705
$ws.Cells.Item(20, 5).Value2 = 0.2967921717790816
$ws.Cells.Item(20, 7).Value2 = 0.002614973350882471
$ws.Cells.Item(20, 9).Value2 = 2.388023175990966
$ws.Cells.Item(20, 10).Value2 = 0.156318211043498
$ws.Cells.Item(20, 11).Value2 = 2.243141198150568
$ws.Cells.Item(20, 14).Value2 = 3.97396587638346
$ws.Cells.Item(21, 2).Value2 = 1.568834220968768
$ws.Cells.Item(21, 3).Value2 = 0.7721060785659688
$ws.Cells.Item(21, 4).Value2 = 0.7375730075112301
$ws.Cells.Item(21, 5).Value2 = 0.304229604398536
$ws.Cells.Item(21, 7).Value2 = 0.002608792925460806
$ws.Cells.Item(21, 9).Value2 = 2.4201530313712
$ws.Cells.Item(21, 10).Value2 = 0.1611954357379801
$ws.Cells.Item(21, 11).Value2 = 2.337933933505894
$ws.Cells.Item(21, 14).Value2 = 3.996574313857309
$ws.Cells.Item(22, 2).Value2 = 1.614229077920811
$ws.Cells.Item(22, 3).Value2 = 0.790013069881752
$ws.Cells.Item(22, 4).Value2 = 0.747904660295859
$ws.Cells.Item(22, 5).Value2 = 0.3092348845617536
$ws.Cells.Item(22, 7).Value2 = 0.002604900448027947
$ws.Cells.Item(22, 9).Value2 = 2.44187783152617
$ws.Cells.Item(22, 10).Value2 = 0.1644636261736565
$ws.Cells.Item(22, 11).Value2 = 2.401074170245124
$ws.Cells.Item(22, 14).Value2 = 4.012021140168173
$ws.Cells.Item(23, 2).Value2 = 1.589926238012367
$ws.Cells.Item(23, 3).Value2 = 0.7804215034852291
$ws.Cells.Item(23, 4).Value2 = 0.7423596145613374
$ws.Cells.Item(23, 5).Value2 = 0.3065502009219472
$ws.Cells.Item(23, 7).Value2 = 0.002606964553123516
$ws.Cells.Item(23, 9).Value2 = 2.430216233133294
$ws.Cells.Item(23, 10).Value2 = 0.1627119151989973
$ws.Cells.Item(23, 11).Value2 = 2.367265864809326
$ws.Cells.Item(23, 14).Value2 = 4.003715234499992
$ws.Cells.Item(24, 2).Value2 = 1.499519814854978
$ws.Cells.Item(24, 3).Value2 = 0.7448444322872092
$ws.Cells.Item(24, 4).Value2 = 0.7220346577737757
$ws.Cells.Item(24, 5).Value2 = 0.2966730388340224
$ws.Cells.Item(24, 7).Value2 = 0.002615076855418878
$ws.Cells.Item(24, 9).Value2 = 2.387510237679464
$ws.Cells.Item(24, 10).Value2 = 0.1562398484449261
$ws.Cells.Item(24, 11).Value2 = 2.241611800869862
$ws.Cells.Item(24, 14).Value2 = 3.973607622722142
$ws.Cells.Item(25, 2).Value2 = 1.405464466742046
$ws.Cells.Item(25, 3).Value2 = 0.7080488481744567
$ws.Cells.Item(25, 4).Value2 = 0.7015446980347235
$ws.Cells.Item(25, 5).Value2 = 0.2866341937220724
$ws.Cells.Item(25, 7).Value2 = 0.002624462354695156
$ws.Cells.Item(25, 9).Value2 = 2.344539528661812
$ws.Cells.Item(25, 10).Value2 = 0.1496011562793811
$ws.Cells.Item(25, 11).Value2 = 2.111123009476643
$ws.Cells.Item(25, 14).Value2 = 3.943986490534201
